# Auto-generated edit script: updates the cryptos price/volume table
# to match the "Updated cryptos list" GitHub Actions commit.
#
# Numeric-looking text values (e.g. "4.80", "0.997") are written with a
# leading apostrophe so Excel keeps them as literal text (preserving
# trailing zeros / exact formatting) instead of silently coercing them
# to floating point numbers, matching the original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.145.54"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "2.740.05"
$ws.Range("E3").Value = "  -5.72%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'508.35"
$ws.Range("E5").Value = "  -3.38%  "
$ws.Range("D6").Value = "'142.84"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'0.535"
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("D9").Value = "2.752.82"
$ws.Range("E9").Value = "  -5.37%  "
$ws.Range("D10").Value = "'6.09"
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").Value = "'0.105"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").Value = "'0.352"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").Value = "3.213.25"
$ws.Range("E14").Value = "  -5.63%  "
$ws.Range("D15").Value = "59.160.07"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").Value = "'21.89"
$ws.Range("E16").Value = "  -3.12%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000137"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.749.68"
$ws.Range("E18").Value = "  -5.24%  "
$ws.Range("D19").Value = "'4.80"
$ws.Range("E19").Value = "  -2.48%  "
$ws.Range("D20").Value = "'11.10"
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("D21").Value = "'347.60"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("D22").Value = "'6.29"
$ws.Range("E22").Value = "  -3.87%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'63.57"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "'0.430"
$ws.Range("E26").Value = "  -3.69%  "
$ws.Range("D27").Value = "'0.174"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("D28").Value = "'0.996"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "0.0₃0847"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "'7.57"
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "'1.62"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").Value = "'19.32"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'149.81"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").Value = "'4.24"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").Value = "'5.45"
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").Value = "'0.963"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").Value = "'1.14"
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("D39").Value = "'36.15"
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("D40").Value = "'1.40"
$ws.Range("E40").Value = "  -4.67%  "
$ws.Range("D41").Value = "'3.56"
$ws.Range("E41").Value = "  -2.05%  "
$ws.Range("D42").Value = "2.193.00"
$ws.Range("E42").Value = "  -5.90%  "
$ws.Range("D43").Value = "'0.0561"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.995"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.609"
$ws.Range("E45").Value = "  -5.50%  "
$ws.Range("D46").Value = "'19.20"
$ws.Range("E46").Value = "  -7.01%  "
$ws.Range("D47").Value = "'4.80"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").Value = "'0.0227"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").Value = "'0.0889"
$ws.Range("E50").Value = "  -3.70%  "
$ws.Range("D51").Value = "'18.26"
$ws.Range("E51").Value = "  +0.58%  "
